$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the "25 Iterations" row (the Edge row of each 3-row block) to
# reflect the new benchmark methodology: 75 total iterations, 25 per
# computing model.
$newLabel = "75 Iterations (25 Each)"
$ws.Range("A4").Value = $newLabel
$ws.Range("A10").Value = $newLabel
$ws.Range("A16").Value = $newLabel
$ws.Range("A22").Value = $newLabel
$ws.Range("A28").Value = $newLabel
$ws.Range("A34").Value = $newLabel
$ws.Range("A40").Value = $newLabel
$ws.Range("A46").Value = $newLabel

# Stray marker cell next to the Large Size - OCR Edge row.
$ws.Range("H34").Value = " "

# Widen the label / timing columns so the new, longer text is readable.
$ws.Columns("A").ColumnWidth = 24.833333333333332
$ws.Columns("D").ColumnWidth = 20.5
$ws.Columns("E").ColumnWidth = 26

# Restore view state: zoom in a bit and leave the selection on D35.
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 132
$ws.Range("D35").Select() | Out-Null
